$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Risks")

$ws.Range("A4").Value = "99aa059c-c144-4966-bbca-95917aa70b10"
$ws.Range("B4").Value = 0.5649999999999999
$ws.Range("C4").Value = 0.131
$ws.Range("D4").Value = 1
$ws.Range("E4").Value = "Mitigation needed"
